$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to remain text so numeric-looking
# strings (e.g. "1.00", "61.429.35") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "61.429.35"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").Value = "3.370.43"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "571.54"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").Value = "136.74"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.367.55"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").Value = "7.46"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("D13").Value = "3.949.65"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("E14").Value = "  +1.89%  "

# Row 15
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("D16").Value = "25.87"
$ws.Range("E16").Value = "  +2.56%  "

# Row 17
$ws.Range("D17").Value = "3.375.62"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("D18").Value = "61.542.53"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "13.94"
$ws.Range("E19").Value = "  +0.31%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "5.88"
$ws.Range("E20").Value = "  -0.13%  "

# Row 21
$ws.Range("D21").Value = "9.32"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("D22").Value = "376.04"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("D23").Value = "0.553"
$ws.Range("E23").Value = "  -2.48%  "

# Row 24
$ws.Range("D24").Value = "3.512.76"
$ws.Range("E24").Value = "  +0.33%  "

# Row 25
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").Value = "71.46"
$ws.Range("E26").Value = "  +1.23%  "

# Row 27
$ws.Range("E27").Value = "  +5.10%  "

# Row 28
$ws.Range("D28").Value = "1.73"
$ws.Range("E28").Value = "  +4.76%  "

# Row 29
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  -2.87%  "

# Row 30
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
$ws.Range("E31").Value = "  +3.10%  "

# Row 32
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("E33").Value = "  +1.40%  "

# Row 34
$ws.Range("E34").Value = "  +0.08%  "

# Row 35
$ws.Range("E35").Value = "  +0.65%  "

# Row 36
$ws.Range("E36").Value = "  -6.99%  "

# Row 37
$ws.Range("E37").Value = "  -3.05%  "

# Row 38
$ws.Range("E38").Value = "  -1.27%  "

# Row 39
$ws.Range("D39").Value = "164.89"
$ws.Range("E39").Value = "  +2.29%  "

# Row 40
$ws.Range("E40").Value = "  -2.50%  "

# Row 41
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.773"
$ws.Range("E42").Value = "  +1.91%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.71"
$ws.Range("E43").Value = "  +1.04%  "

# Row 44
$ws.Range("D44").Value = "41.56"
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("E45").Value = "  +0.69%  "

# Row 46
$ws.Range("D46").Value = "4.37"
$ws.Range("E46").Value = "  -1.18%  "

# Row 47
$ws.Range("D47").Value = "24.64"
$ws.Range("E47").Value = "  +6.43%  "

# Row 48
$ws.Range("D48").Value = "6.84"
$ws.Range("E48").Value = "  -1.80%  "

# Row 49
$ws.Range("D49").Value = "22.73"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("D50").Value = "2.353.04"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51
$ws.Range("E51").Value = "  +0.51%  "
